# Insert a new data row at row 93 (pushing existing rows 93-133 down to 94-134)
# and populate it with the new weekly Achicoria price record for
# "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 93, shifting rows 93:133 to 94:134
$ws.Rows.Item(93).Insert()

# Fill in the new row 93 with the new record's values
$ws.Range("A93").Value = 10
$ws.Range("B93").Value = "Vega Modelo de Temuco"
$ws.Range("C93").Value = "La Araucanía"
$ws.Range("D93").Value = 45141
$ws.Range("E93").Value = 9
$ws.Range("F93").Value = 100112010
$ws.Range("G93").Value = "Achicoria"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 300
$ws.Range("K93").Value = 10000
$ws.Range("L93").Value = 10000
$ws.Range("M93").Value = 10000
$ws.Range("N93").Value = "$/caja 18 unidades"
$ws.Range("O93").Value = "Región Metropolitana"
$ws.Range("P93").Value = 556
$ws.Range("Q93").Value = 18
$ws.Range("R93").Value = "Hortaliza"

# Ensure date column D keeps the same date style as the rest of column D
$ws.Range("D93").NumberFormat = $ws.Range("D94").NumberFormat
